# LOM3043.xlsx update
# - Inserts a new row at position 13 (shifting the existing "Docentes
#   responsáveis:" value row and everything below it down by one row),
#   then fixes up the text content so every label in column A lines up
#   with its correct description in columns B/C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 13. Excel shifts rows 13:25 down to 14:26 and
# carries row heights / styles with them automatically.
$ws.Rows(13).Insert()

# The insert cloned column A's style into the new A13; the target layout
# has no cell there at all, so drop it.
$ws.Range("A13").Clear()

# Row 13 (new): "Docentes responsáveis:" value, moved up from the old
# row 10 position. Copy the body-text formatting (wrap text, not bold /
# red text) from the row below before filling it in.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B13").Value = "5840622 - Miguel Justino Ribeiro Barboza"
$ws.Range("C13").Value = "5840622 - Miguel Justino Ribeiro Barboza"
$excel.CutCopyMode = $false

# Row 10: "Objetivos:" body text.
$ws.Range("B10").Value = "Apresentar os princípios básicos da Seleção de Materiais para aplicação em Engenharia."
$ws.Range("C10").Value = "Apresentar os princípios básicos da Seleção de Materiais para aplicação em Engenharia."

# Row 14: "Programa resumido:" body text.
$ws.Range("B14").Value = "Aspectos gerais e critérios de seleção de materiais estruturais. Aspectos dos principais mecanismos de falha em componentes estruturais. Seleção de materiais e análise para diferentes modos de carregamento. Seleção de materiais sob diferentes condições de temperatura. Materiais resistentes à corrosão e oxidação. Tribologia: atrito e desgaste. Tratamentos superficiais."
$ws.Range("C14").Value = "Aspectos gerais e critérios de seleção de materiais estruturais. Aspectos dos principais mecanismos de falha em componentes estruturais. Seleção de materiais e análise para diferentes modos de carregamento. Seleção de materiais sob diferentes condições de temperatura. Materiais resistentes à corrosão e oxidação. Tribologia: atrito e desgaste. Tratamentos superficiais."

# Row 16: "Programa:" body text.
$ws.Range("B16").Value = "1. Principais mecanismos de falha em componentes estruturais: efeitos do meio e temperatura. Critérios de falha. 2. Seleção de materiais para aplicações sob a ação de cargas estáticas. Materiais metálicos, cerâmicos, poliméricos e compósitos. 3. Seleção de materiais para aplicações sob a ação de cargas dinâmicas: O fenômeno da fadiga e efeitos da presença de entalhes em componentes mecânicos. 4. Critérios de seleção de materiais para aplicações em temperaturas elevadas. O fenômeno da fluência e a tolerância ao dano. Seleção de materiais para alta temperatura. Aços especiais, superligas, materiais cerâmicos e compósitos. 5. Materiais para temperaturas criogênicas. A transição dúctil-frágil. 6. Aspectos fundamentais do estudo de tribologia: desgaste, atrito e tratamentos superficiais. 7. Fundamentos, seleção e proteção contra oxidação. 8. Seleção de materiais em meios corrosivos. Corrosão sob tensão."
$ws.Range("C16").Value = "1. Principais mecanismos de falha em componentes estruturais: efeitos do meio e temperatura. Critérios de falha. 2. Seleção de materiais para aplicações sob a ação de cargas estáticas. Materiais metálicos, cerâmicos, poliméricos e compósitos. 3. Seleção de materiais para aplicações sob a ação de cargas dinâmicas: O fenômeno da fadiga e efeitos da presença de entalhes em componentes mecânicos. 4. Critérios de seleção de materiais para aplicações em temperaturas elevadas. O fenômeno da fluência e a tolerância ao dano. Seleção de materiais para alta temperatura. Aços especiais, superligas, materiais cerâmicos e compósitos. 5. Materiais para temperaturas criogênicas. A transição dúctil-frágil. 6. Aspectos fundamentais do estudo de tribologia: desgaste, atrito e tratamentos superficiais. 7. Fundamentos, seleção e proteção contra oxidação. 8. Seleção de materiais em meios corrosivos. Corrosão sob tensão."

# Row 19: "Método:" body text.
$ws.Range("B19").Value = "Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa."
$ws.Range("C19").Value = "Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa."

# Row 20: "Critério:" body text.
$ws.Range("B20").Value = "A média do semestre será computada com base na relação:M=(P1+2P2)/3"
$ws.Range("C20").Value = "A média do semestre será computada com base na relação:M=(P1+2P2)/3"

# Row 21: "Norma de recuperação:" body text.
$ws.Range("B21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2"
$ws.Range("C21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2"

# Row 22: "Bibliografia:" body text.
$ws.Range("B22").Value = "1.Ashby, M. F. Materials Selection in Mechanical Design, Butterworth, Oxford, 2005. 2. ASM Metals Handbook - Properties and Selection: Irons, Steels and High - Performance Alloys - v.1 - 1990. 3. ASM Metals Handbook - Properties and Selection: Nonferrous Alloys and Special - Purpose Materials - v.2 - 1990. 4. Meyers, M.; Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009. 5. Van Vlack, L.H., Propriedades dos Materiais Cerâmicos. Ed. Edgard Blücher Ltda., 1973. 6. Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall, 1999. 7. Biasotto, E., Polímeros como Materiais de Engenharia. Ed. Edgard Blücher Ltda., 1991. 8. Rosen, S.L., Fundamental Principles of Polymeric Materials. Ed. John Wiley & Sons, Inc., 1993. 9. Bhushan, B. Introduction to Tribology, 2nd Edition, John Wiley & Sons. 2013. 10. Roberge, P. R. Corrosion engineering: principles and practice. The McGraw-Hill Companies, Inc., 2008. 11. Gentil, V. Corrosão, Ed. LTC, 2011. 12. Crane, F.A., Charles, J.A., Selection of Engineering Materials, Butterworth, 1984. 13. Chiaverini, V., Aços e Ferros Fundidos, Associação Brasileira de Materiais - ABM, São Paulo, 1988. 14. Reed, R. C. The superalloys: fundamentals and applications. Ed. Cambridge, USA, 2006."
$ws.Range("C22").Value = "1.Ashby, M. F. Materials Selection in Mechanical Design, Butterworth, Oxford, 2005. 2. ASM Metals Handbook - Properties and Selection: Irons, Steels and High - Performance Alloys - v.1 - 1990. 3. ASM Metals Handbook - Properties and Selection: Nonferrous Alloys and Special - Purpose Materials - v.2 - 1990. 4. Meyers, M.; Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009. 5. Van Vlack, L.H., Propriedades dos Materiais Cerâmicos. Ed. Edgard Blücher Ltda., 1973. 6. Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall, 1999. 7. Biasotto, E., Polímeros como Materiais de Engenharia. Ed. Edgard Blücher Ltda., 1991. 8. Rosen, S.L., Fundamental Principles of Polymeric Materials. Ed. John Wiley & Sons, Inc., 1993. 9. Bhushan, B. Introduction to Tribology, 2nd Edition, John Wiley & Sons. 2013. 10. Roberge, P. R. Corrosion engineering: principles and practice. The McGraw-Hill Companies, Inc., 2008. 11. Gentil, V. Corrosão, Ed. LTC, 2011. 12. Crane, F.A., Charles, J.A., Selection of Engineering Materials, Butterworth, 1984. 13. Chiaverini, V., Aços e Ferros Fundidos, Associação Brasileira de Materiais - ABM, São Paulo, 1988. 14. Reed, R. C. The superalloys: fundamentals and applications. Ed. Cambridge, USA, 2006."

Write-Output "Applied LOM3043 content updates"
